# Sex = null -> sex = "u"
# The "dimensions" worksheet holds rows describing stat-graph dimensions in
# columns A (computed key), B/C (dimension names) and D (sex: m / f / null).
# D4 is the only literal "sex" cell holding the sentinel "null"; D7, D10 and
# D13 are shared formulas (=D4, =D7, =D10 respectively) that just mirror it,
# and column A as well as the whole "dimensions stats" sheet are formula
# driven off column D. So updating the single literal cell D4 is enough for
# everything downstream (A4/A7/A10/A13, D7/D10/D13, and sheet "dimensions
# stats") to recalculate to the new "u" sentinel automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("dimensions")

$ws.Range("D4").Value = "u"

$ws.Activate()
$ws.Range("D5").Select()
